# Applies the "Updated cryptos list" data refresh to sheet1 (Price / Volume(1h) columns).
# For cells whose new text looks like a plain number (e.g. "1.004"), we force the cell to
# stay text (NumberFormat "@") and then restore the default "Normal" style so no stray
# formatting is left behind once Excel has re-parsed the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.970.67'
$ws.Range("E2").Value = '  +0.36%  '
# Row 3
$ws.Range("D3").Value = '1.924.54'
$ws.Range("E3").Value = '  +1.31%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4587'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.15%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3819'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.26%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07750'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.29%  '
# Row 10
$ws.Range("E10").Value = '  -0.03%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.57'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.54%  '
# Row 12
$ws.Range("D12").Value = '1.951.10'
$ws.Range("E12").Value = '  +1.16%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.708'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.67%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.972'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.21%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06985'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.91%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.96%  '
# Row 17
$ws.Range("E17").Value = '  +0.00%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009491'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.44%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.20%  '
# Row 20
$ws.Range("E20").Value = '  -0.07%  '
# Row 21
$ws.Range("D21").Value = '28.987.93'
$ws.Range("E21").Value = '  +0.52%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.345'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.17%  '
# Row 23
$ws.Range("E23").Value = '  +1.67%  '
# Row 24
$ws.Range("D24").Value = '2.136.87'
$ws.Range("E24").Value = '  -0.73%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.59%  '
# Row 27
$ws.Range("E27").Value = '  -0.68%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.622'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.82%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
# Row 30
$ws.Range("E30").Value = '  +0.30%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09315'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.51%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8653'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.46%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.110'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.23%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.246'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.55%  '
# Row 35
$ws.Range("E35").Value = '  -0.13%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05694'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.04%  '
# Row 37
$ws.Range("E37").Value = '  +0.49%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.004'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.17%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02053'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.97%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.104'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.58%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.465'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.01%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5508'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.05%  '
# Row 43
$ws.Range("E43").Value = '  +0.15%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.343'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.42%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002792'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +17.59%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.184'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.25%  '
# Row 47
$ws.Range("E47").Value = '  -0.33%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06942'
$ws.Range("D48").Style = "Normal"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.21'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.21%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.46%  '
# Row 51
$ws.Range("E51").Value = '  -0.23%  '
